# Weekly update: insert two new price records (rows 198-199) for
# "Poroto granado" at Vega Central Mapocho de Santiago, shifting the
# existing historical rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the current row 198, pushing rows 198:252
# down to 200:254 (dimension grows from A1:R252 to A1:R254).
$ws.Rows("198:199").Insert()

# New row 198
$ws.Range("A198").Value = 9
$ws.Range("B198").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C198").Value = "Metropolitana"
$ws.Range("D198").Value = 44627
$ws.Range("E198").Value = 13
$ws.Range("F198").Value = 100112030
$ws.Range("G198").Value = "Poroto granado"
$ws.Range("H198").Value = "Sin especificar"
$ws.Range("I198").Value = "Primera"
$ws.Range("J198").Value = 52
$ws.Range("K198").Value = 25000
$ws.Range("L198").Value = 27000
$ws.Range("M198").Value = 26000
$ws.Range("N198").Value = "`$/saco 25 kilos"
$ws.Range("O198").Value = "Región Metropolitana"
$ws.Range("P198").Value = 1040
$ws.Range("Q198").Value = 25
$ws.Range("R198").Value = "Hortaliza"

# New row 199
$ws.Range("A199").Value = 9
$ws.Range("B199").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C199").Value = "Metropolitana"
$ws.Range("D199").Value = 44627
$ws.Range("E199").Value = 13
$ws.Range("F199").Value = 100112030
$ws.Range("G199").Value = "Poroto granado"
$ws.Range("H199").Value = "Sin especificar"
$ws.Range("I199").Value = "Primera"
$ws.Range("J199").Value = 43
$ws.Range("K199").Value = 25000
$ws.Range("L199").Value = 27000
$ws.Range("M199").Value = 26023
$ws.Range("N199").Value = "`$/saco 25 kilos"
$ws.Range("O199").Value = "Región de O'Higgins"
$ws.Range("P199").Value = 1041
$ws.Range("Q199").Value = 25
$ws.Range("R199").Value = "Hortaliza"
